$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores figures as plain text in the source data (inline
# strings), e.g. "1.400" / "215.77". A handful of the new values round-trip
# through the General number parser, which would silently convert them to
# numbers (and mangle trailing zeros, e.g. "1.400" -> 1.4). Force those
# specific cells to Text format first so they stay text, matching the source.
$textCells = @("D5", "D9", "D10", "D11", "D15", "D17", "D21", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.073.94'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '1.647.10'
$ws.Range("E3").Value = '  -1.34%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '215.77'
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").Value = '0.06372'
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("D10").Value = '20.89'
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").Value = '0.07673'
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("D12").Value = '1.647.15'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").Value = '1.869.80'
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("D15").Value = '0.5552'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").Value = '0.0₅8332'
$ws.Range("E16").Value = '  +3.85%  '
$ws.Range("D17").Value = '64.98'
$ws.Range("D18").Value = '26.073.68'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '188.71'
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").Value = '6.277'
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").Value = '146.04'
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("D27").Value = '7.427'
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").Value = '15.86'
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").Value = '1.400'
$ws.Range("E29").Value = '  +3.08%  '
$ws.Range("D30").Value = '0.05972'
$ws.Range("E30").Value = '  -5.13%  '
$ws.Range("D31").Value = '1.269'
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D32").Value = '3.404'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("D33").Value = '3.412'
$ws.Range("E33").Value = '  -3.26%  '
$ws.Range("D34").Value = '1.659'
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").Value = '0.9987'
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").Value = '2.753'
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("E38").Value = '  -6.03%  '
$ws.Range("D39").Value = '0.01611'
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").Value = '0.8583'
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("D41").Value = '5.841'
$ws.Range("E41").Value = '  -3.72%  '
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '1.029.51'
$ws.Range("E43").Value = '  -7.76%  '
$ws.Range("D44").Value = '98.83'
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("D45").Value = '1.796.79'
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("D47").Value = '55.87'
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").Value = '1.005'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '8.076'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '0.05154'
$ws.Range("D51").Value = '0.4216'
$ws.Range("E51").Value = '  -0.53%  '
